$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 4's formatting down into the new row 5, then overwrite
# with the new MobileNetV2 experiment's values.
$ws.Range("A4:H4").Copy() | Out-Null
$ws.Range("A5").Select() | Out-Null
$ws.Paste() | Out-Null

$ws.Range("A5").Value = 220609
$ws.Range("B5").Value = "박영서"
$ws.Range("C5").Value = "MobileNetV2"
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 0.8664
$ws.Range("F5").Value = 0.8575
$ws.Range("G5").Value = 0.5156
$ws.Range("H5").Value = 0.5129

# Widen column C (model name) now that "MobileNetV2" no longer fits the
# bestFit width, and drop the stale bestFit flag in favor of a fixed width.
$ws.Columns.Item(3).ColumnWidth = 11.6640625

# Leave the selection on the newly added accuracy cells.
$ws.Range("G5:H5").Select() | Out-Null
